$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the trailing break-out rows (5-7) and the now-unused Supplier Name
# column E; the remaining four columns collapse into a single flat table.
$ws.Rows("5:7").Delete()
$ws.Columns("E:E").Delete()
$ws.Cells.UnMerge()

# New header row
$ws.Range("A1").Value = "SR"
$ws.Range("B1").Value = "Part Number"
$ws.Range("C1").Value = "Supplier Code"
$ws.Range("D1").Value = "Supplier Name"

# New data rows
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "num1"
$ws.Range("C2").Value = "S111"
$ws.Range("D2").Value = "ABCg"

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "num2"
$ws.Range("C3").Value = "S112"
$ws.Range("D3").Value = "SBNk"

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "num3"
$ws.Range("C4").Value = "S113"
$ws.Range("D4").Value = "SBNs"

# Supplier Code / Part Number columns in the data rows are plain bordered
# cells (no center alignment) now, unlike the old SR/Part-Number columns.
$ws.Range("C2:C4").HorizontalAlignment = -4131
$ws.Range("C2:C4").VerticalAlignment = -4131

# Last cell (D4) reverts to the default, unbordered style.
$ws.Range("D4").Borders.LineStyle = -4142
$ws.Range("D4").Font.Bold = $false

# Sheet view / selection bookkeeping
$ws.Range("D4").Select()

$wb.Windows.Item(1).WindowState = -4143
